$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3622.8635
$ws.Range("J69").Value = 4357.3
$ws.Range("L69").Value = 13071.9
$ws.Range("N69").Value = -14819.9
$ws.Range("H72").Value = 3622.8635
$ws.Range("J72").Value = 4357.3
$ws.Range("L72").Value = 39215.7
$ws.Range("N72").Value = -47951.7
$ws.Range("H80").Value = 881.4
$ws.Range("I80").Value = 292.6
$ws.Range("J80").Value = 1470.2
$ws.Range("K80").Value = 877.8000000000001
$ws.Range("L80").Value = 4410.6
$ws.Range("M80").Value = 120.1999999999999
$ws.Range("N80").Value = -6406.6
$ws.Range("H83").Value = 881.4
$ws.Range("I83").Value = 292.6
$ws.Range("J83").Value = 1470.2
$ws.Range("K83").Value = 2633.4
$ws.Range("L83").Value = 13231.8
$ws.Range("M83").Value = 2358.6
$ws.Range("N83").Value = -23215.8
$ws.Range("H113").Value = 8899.923000000001
$ws.Range("I113").Value = 3599.8572
$ws.Range("J113").Value = 10852.579
$ws.Range("K113").Value = 3599.8572
$ws.Range("L113").Value = 10852.579
$ws.Range("M113").Value = -345.8571999999999
$ws.Range("N113").Value = -17360.579
$ws.Range("H129").Value = 1016.3137
$ws.Range("J129").Value = 1096.1818
$ws.Range("L129").Value = 3288.5454
$ws.Range("N129").Value = -13288.5454
$ws.Range("H132").Value = 35549.11
$ws.Range("I132").Value = 45754.18
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 137262.54
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -134732.54
$ws.Range("N132").Value = -16460
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 7000
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 21000
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -18450
$ws.Range("N137").Value = -8100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 675
$ws.Range("I4").Value = 733.3333
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 733.3333
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -617.3333
$ws.Range("N4").Value = -732
$ws.Range("H32").Value = 6520.828
$ws.Range("I32").Value = 5162.788
$ws.Range("J32").Value = 20950
$ws.Range("K32").Value = 5162.788
$ws.Range("L32").Value = 20950
$ws.Range("M32").Value = -4875.788
$ws.Range("N32").Value = -21524
$ws.Range("H45").Value = 1207.2727
$ws.Range("I45").Value = 1128
$ws.Range("K45").Value = 1128
$ws.Range("M45").Value = -751
$ws.Range("H112").Value = 15387
$ws.Range("J112").Value = 15387
$ws.Range("L112").Value = 15387
$ws.Range("N112").Value = -18341
$ws.Range("H139").Value = 46236.668
$ws.Range("J139").Value = 46236.668
$ws.Range("L139").Value = 46236.668
$ws.Range("N139").Value = -56516.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H140").Value = 54427.145
$ws.Range("J140").Value = 54427.145
$ws.Range("L140").Value = 54427.145
$ws.Range("N140").Value = -64787.145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3793.0444
$ws.Range("I31").Value = 1538.8704
$ws.Range("J31").Value = 7174.3057
$ws.Range("K31").Value = 1538.8704
$ws.Range("L31").Value = 7174.3057
$ws.Range("M31").Value = -1243.8704
$ws.Range("N31").Value = -7764.3057
$ws.Range("H34").Value = 3793.0444
$ws.Range("I34").Value = 1538.8704
$ws.Range("J34").Value = 7174.3057
$ws.Range("K34").Value = 1538.8704
$ws.Range("L34").Value = 7174.3057
$ws.Range("M34").Value = -1336.8704
$ws.Range("N34").Value = -7578.3057
$ws.Range("H107").Value = 380.83334
$ws.Range("I107").Value = 308.33334
$ws.Range("J107").Value = 453.33334
$ws.Range("K107").Value = 308.33334
$ws.Range("L107").Value = 453.33334
$ws.Range("M107").Value = 1611.66666
$ws.Range("N107").Value = -4293.33334
$ws.Range("H138").Value = 38552.5
$ws.Range("J138").Value = 38552.5
$ws.Range("L138").Value = 38552.5
$ws.Range("N138").Value = -48832.5
$ws.Range("H140").Value = 54783.168
$ws.Range("J140").Value = 54783.168
$ws.Range("L140").Value = 54783.168
$ws.Range("N140").Value = -65143.168

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 660.1429000000001
$ws.Range("I15").Value = 30.25
$ws.Range("J15").Value = 1500
$ws.Range("K15").Value = 90.75
$ws.Range("L15").Value = 4500
$ws.Range("M15").Value = 49.25
$ws.Range("N15").Value = -4780
$ws.Range("H122").Value = 969.7143
$ws.Range("I122").Value = 520.34784
$ws.Range("J122").Value = 3036.8
$ws.Range("K122").Value = 4683.130560000001
$ws.Range("L122").Value = 27331.2
$ws.Range("M122").Value = -2233.130560000001
$ws.Range("N122").Value = -32231.2
$ws.Range("H138").Value = 1359.6923
$ws.Range("I138").Value = 959.5625
$ws.Range("J138").Value = 1999.9
$ws.Range("K138").Value = 2878.6875
$ws.Range("L138").Value = 5999.700000000001
$ws.Range("M138").Value = 2261.3125
$ws.Range("N138").Value = -16279.7

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2362.4243
$ws.Range("I126").Value = 2250.2222
$ws.Range("J126").Value = 2497.0667
$ws.Range("K126").Value = 6750.6666
$ws.Range("L126").Value = 7491.2001
$ws.Range("M126").Value = -4280.6666
$ws.Range("N126").Value = -12431.2001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 673.1818
$ws.Range("I22").Value = 630
$ws.Range("J22").Value = 703.0769
$ws.Range("K22").Value = 630
$ws.Range("L22").Value = 703.0769
$ws.Range("M22").Value = -335
$ws.Range("N22").Value = -1293.0769
$ws.Range("H27").Value = 673.1818
$ws.Range("I27").Value = 630
$ws.Range("J27").Value = 703.0769
$ws.Range("K27").Value = 630
$ws.Range("L27").Value = 703.0769
$ws.Range("M27").Value = -523
$ws.Range("N27").Value = -917.0769
$ws.Range("H40").Value = 49890.91
$ws.Range("I40").Value = 66750
$ws.Range("J40").Value = 4933.3335
$ws.Range("K40").Value = 66750
$ws.Range("L40").Value = 4933.3335
$ws.Range("M40").Value = -66614
$ws.Range("N40").Value = -5205.3335
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H133").Value = 43454
$ws.Range("J133").Value = 43454
$ws.Range("L133").Value = 43454
$ws.Range("N133").Value = -48514
$ws.Range("H136").Value = 10103633
$ws.Range("I136").Value = 2611.0833
$ws.Range("K136").Value = 7833.249899999999
$ws.Range("M136").Value = -5283.249899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2107.4666
$ws.Range("I132").Value = 1426.7142
$ws.Range("J132").Value = 3695.889
$ws.Range("K132").Value = 4280.142599999999
$ws.Range("L132").Value = 11087.667
$ws.Range("M132").Value = -1750.142599999999
$ws.Range("N132").Value = -16147.667
